# Fixed Stimulus Absolute Timestamps
# Renames each task-order sheet (new timestamp suffix) and updates the
# stimulus-file CSV names / order referenced on each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504778286298823"
$ws1.Range("B2").Value = "go_stims-16504778285868807.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778286118772.csv"
$ws1.Range("B4").Value = "go_stims-165047782861388.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778286278772.csv"

# --- Sheet 2: NB_TO ----------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778303419137"
$ws2.Range("B2").Value = "ZB-match_6-16504778289208817.csv"
$ws2.Range("B3").Value = "TB-16504778303239126.csv"
$ws2.Range("B4").Value = "OB-16504778301508794.csv"
$ws2.Range("B5").Value = "OB-1650477829524913.csv"
$ws2.Range("B6").Value = "TB-16504778302879121.csv"
$ws2.Range("B7").Value = "TB-16504778302548795.csv"
$ws2.Range("B8").Value = "ZB-match_1-1650477829000878.csv"
$ws2.Range("B9").Value = "OB-1650477830040911.csv"
$ws2.Range("B10").Value = "ZB-match_0-1650477829122879.csv"

# --- Sheet 3: RS_TO (only the name changes) -----------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16504778303438828"

# --- Sheet 4: TOL_TO -----------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778303908784"
$ws4.Range("B2").Value = "MM_stims-1650477830357895.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778303458781.csv"
$ws4.Range("B4").Value = "MM_stims-16504778303739119.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477830358881.csv"
$ws4.Range("B6").Value = "MM_stims-1650477830389907.csv"
$ws4.Range("B7").Value = "ZM_stims-165047783037488.csv"

# --- Sheet 5: vSAT_TO ----------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778304538765"
$ws5.Range("B2").Value = "SAT_stims-1650477830393881.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778304058962.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778304219108.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650477830437911.csv"
